$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert the four new header paragraphs (blank / "Aluno: ..." / "Materia:
#    ..." / blank) in front of the existing first paragraph ("CRISE DE
#    SOFTWARE:"). We build the exact OOXML for the new paragraphs and insert
#    it via Range.InsertXML at a collapsed range at the very start of the
#    story so the existing content is untouched and simply pushed down.
# ---------------------------------------------------------------------------

$newParagraphsXml = @'
<w:p><w:pPr><w:pStyle w:val="Normal"/><w:jc w:val="center"/><w:rPr><w:b/><w:b/><w:bCs/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr/></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Normal"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t>Aluno: Henrique S. Paini</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Normal"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t>Mat</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Calibri" w:cs="" w:cstheme="minorBidi" w:eastAsiaTheme="minorHAnsi"/><w:color w:val="auto"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="pt-BR" w:eastAsia="en-US" w:bidi="ar-SA"/></w:rPr><w:t>\u00e9ria: Engenharia de Software I.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Normal"/><w:jc w:val="center"/><w:rPr><w:b/><w:b/><w:bCs/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr/></w:r></w:p>
'@

# PowerShell here-strings don't interpret \u escapes; substitute the
# accented character explicitly so the literal text stays exactly
# "Materia: Engenharia de Software I." (with the accented e).
$newParagraphsXml = $newParagraphsXml.Replace('\u00e9', [string]([char]0x00E9))

$package = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' + $newParagraphsXml + '</w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint = $d.Range(0, 0)
$insertionPoint.InsertXML($package)

# ---------------------------------------------------------------------------
# 2. Flag the "Normal" style's paragraph properties so automatic hyphenation
#    is suppressed (adds <w:suppressAutoHyphens/> right after
#    <w:widowControl/> in word/styles.xml).
# ---------------------------------------------------------------------------

$normalStyle = $d.Styles("Normal")
$normalStyle.ParagraphFormat.Hyphenation = $false
